$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values (27CON/04CON/27v04 -> F1CON/F2CON/F1vF2)
$ws.Range("A2").Value = "F1CON"
$ws.Range("B2").Value = "F2CON"
$ws.Range("C2").Value = "F1vF2"

# Add new comparison rows 3-7
$ws.Range("A3").Value = "F1CON"
$ws.Range("B3").Value = "F3CON"
$ws.Range("C3").Value = "F1vF3"

$ws.Range("A4").Value = "F2CON"
$ws.Range("B4").Value = "F1CON"
$ws.Range("C4").Value = "F2vF1"

$ws.Range("A5").Value = "F2CON"
$ws.Range("B5").Value = "F3CON"
$ws.Range("C5").Value = "F2vF3"

$ws.Range("A6").Value = "F3CON"
$ws.Range("B6").Value = "F1CON"
$ws.Range("C6").Value = "F3vF1"

$ws.Range("A7").Value = "F3CON"
$ws.Range("B7").Value = "F2CON"
$ws.Range("C7").Value = "F3vF2"

# Update selection to match target (D5)
$ws.Range("D5").Select()
